$wb = $excel.ActiveWorkbook

# --- Shared strings: add new label used by sheet1 row 7 ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Sheet1 ("All_model_short"): update regression re-run values
$ws1.Range("A7").Value = "DistCenter_res_pc"
$ws1.Range("B2").Value = -1615.7011365527501
$ws1.Range("C2").Value = [double]"2.3150861691328E-7"
$ws1.Range("B3").Value = -2096.41471530015
$ws1.Range("C3").Value = [double]"2.9379876256380003E-11"
$ws1.Range("B4").Value = -1365.70154334427
$ws1.Range("C4").Value = [double]"3.0171061763258E-5"
$ws1.Range("B5").Value = 19.660293114867901
$ws1.Range("C5").Value = [double]"5.1902917160977996E-3"
$ws1.Range("B6").Value = 73.971545363940507
$ws1.Range("C6").Value = [double]"1.18835537232273E-26"
$ws1.Range("B7").Value = 9.0211740009986805
$ws1.Range("C7").Value = [double]"1.50492399107082E-35"
$ws1.Range("B8").Value = 60.965563012219697
$ws1.Range("C8").Value = [double]"9.70869484935E-4"
$ws1.Range("B9").Value = -3.4808134458465898
$ws1.Range("C9").Value = [double]"1.5234106167741399E-6"
$ws1.Range("B10").Value = -3.5778093292700799
$ws1.Range("C10").Value = 0.36689261745722501
$ws1.Range("B11").Value = -17.9553753129216
$ws1.Range("C11").Value = [double]"8.8965345328113105E-12"
$ws1.Range("B12").Value = 1.8430814938706399
$ws1.Range("C12").Value = 0.371888426652624
$ws1.Range("B13").Value = -4.8266466803171602
$ws1.Range("C13").Value = [double]"6.2999910862285596E-2"
$ws1.Range("B14").Value = -18.7013217205748
$ws1.Range("C14").Value = [double]"2.5177670029870298E-7"

# Sheet2 ("All_model_short_table"): same re-run values (label in A7 unchanged)
$ws2.Range("B2").Value = -1615.7011365527501
$ws2.Range("C2").Value = [double]"2.3150861691328E-7"
$ws2.Range("B3").Value = -2096.41471530015
$ws2.Range("C3").Value = [double]"2.9379876256380003E-11"
$ws2.Range("B4").Value = -1365.70154334427
$ws2.Range("C4").Value = [double]"3.0171061763258E-5"
$ws2.Range("B5").Value = 19.660293114867901
$ws2.Range("C5").Value = [double]"5.1902917160977996E-3"
$ws2.Range("B6").Value = 73.971545363940507
$ws2.Range("C6").Value = [double]"1.18835537232273E-26"
$ws2.Range("B7").Value = 9.0211740009986805
$ws2.Range("C7").Value = [double]"1.50492399107082E-35"
$ws2.Range("B8").Value = 60.965563012219697
$ws2.Range("C8").Value = [double]"9.70869484935E-4"
$ws2.Range("B9").Value = -3.4808134458465898
$ws2.Range("C9").Value = [double]"1.5234106167741399E-6"
$ws2.Range("B10").Value = -3.5778093292700799
$ws2.Range("C10").Value = 0.36689261745722501
$ws2.Range("B11").Value = -17.9553753129216
$ws2.Range("C11").Value = [double]"8.8965345328113105E-12"
$ws2.Range("B12").Value = 1.8430814938706399
$ws2.Range("C12").Value = 0.371888426652624
$ws2.Range("B13").Value = -4.8266466803171602
$ws2.Range("C13").Value = [double]"6.2999910862285596E-2"
$ws2.Range("B14").Value = -18.7013217205748
$ws2.Range("C14").Value = [double]"2.5177670029870298E-7"

# --- Selections / view state ---
$ws1.Range("A2:C14").Select()

$ws2.Range("F13").Select()

# --- Column A width on sheet2 (best-fit-like) ---
$ws2.Columns.Item(1).ColumnWidth = 17
